# Daily attendance processing - 2026-01-16 15:07:22
#
# The "Recorded By" column (G) lists the accounts that touched a given
# attendance session, as a comma-separated string. This run normalises the
# ordering so the human account (dnasr281@gmail.com) is listed first,
# ahead of automated / admin accounts, for every two-entry list affected
# by today's processing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "Recorded By" (column G) whose two comma-separated entries
# need to be swapped so dnasr281@gmail.com appears first.
$rows = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G" + $r)
    $current = $cell.Value2
    $parts = $current -split ", "
    if ($parts.Length -eq 2 -and $parts[1] -eq "dnasr281@gmail.com") {
        $cell.Value = $parts[1] + ", " + $parts[0]
    }
}
